# Keys_Onboarding / TestData.xlsx — "Finish create new propety and verity"
#
# 1. PropertyDetails!M2 (FilePath) picks up the finished local screenshot
#    path used for verification.
# 2. A brand new "TenantDetails" sheet is added (end of the workbook) that
#    captures the tenant used to verify the onboarding flow.

$wb = $excel.ActiveWorkbook

# --- 1. Update the property photo FilePath on PropertyDetails -------------
$propertySheet = $wb.Worksheets.Item("PropertyDetails")
$propertySheet.Range("M2").Value = "C:\Users\Mogan\source\repos\ICOnboardingTask\Keys_Onboarding\Photos\01.jpeg"

# --- 2. Add the new TenantDetails sheet at the end of the workbook --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tenantSheet = $wb.Worksheets.Add($null, $lastSheet)
$tenantSheet.Name = "TenantDetails"

$tenantSheet.Cells.Item(1,1).Value = "TenantEmail"
$tenantSheet.Cells.Item(1,2).Value = "IsMainTenant"
$tenantSheet.Cells.Item(1,3).Value = "FirstName"
$tenantSheet.Cells.Item(1,4).Value = "LastName"
$tenantSheet.Cells.Item(1,5).Value = "Duraion"
$tenantSheet.Cells.Item(1,6).Value = "RentAmount"
$tenantSheet.Cells.Item(1,7).Value = "Paymentfrequency"
$tenantSheet.Cells.Item(1,8).Value = "PaymentDueDay"

$tenantSheet.Cells.Item(2,1).Value = "test@test.com"
$tenantSheet.Hyperlinks.Add($tenantSheet.Range("A2"), "mailto:test@test.com")
$tenantSheet.Cells.Item(2,2).Value = "Yes"
$tenantSheet.Cells.Item(2,3).Value = "test"
$tenantSheet.Cells.Item(2,4).Value = "test"
$tenantSheet.Cells.Item(2,5).Value = 125
$tenantSheet.Cells.Item(2,6).Value = 23
$tenantSheet.Cells.Item(2,7).Value = "Fortnightly"
$tenantSheet.Cells.Item(2,8).Value = 5

# --- 3. Leave PropertyDetails as the active tab ----------------------------
$propertySheet.Activate()
